{"js": "// Fix: summary texts highlights\n// 1) \", e na aplica\u00e7\u00e3o dos paradigmas de \" + \"programa\u00e7\u00e3o funcional\" (bold)\n//    -> \", e na aplica\u00e7\u00e3o dos paradigmas de programa\u00e7\u00e3o \" + \"funcional\" (bold)\n// 2) \"Ex-coordenadora da comunidade sem fins lucrativos \"\n//    -> \"Ex-coordenadora e atual membra da comunidade sem fins lucrativos \"\n\n// --- Edit 1: move \"programa\u00e7\u00e3o \" out of the bold run into the preceding\n// non-bold run, leaving only \"funcional\" bold. ---\nconst boldPhrase = context.document.body.search(\"programa\u00e7\u00e3o funcional\", { matchCase: true });\nboldPhrase.load(\"text\");\nawait context.sync();\n\nif (boldPhrase.items.length > 0) {\n  // Shrink the bold run down to just \"funcional\".\n  boldPhrase.items[0].insertText(\"funcional\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst prefixPhrase = context.document.body.search(\n  \", e na aplica\u00e7\u00e3o dos paradigmas de \",\n  { matchCase: true }\n);\nprefixPhrase.load(\"text\");\nawait context.sync();\n\nif (prefixPhrase.items.length > 0) {\n  // Extend the preceding (non-bold) run so it now ends with \"programa\u00e7\u00e3o \".\n  prefixPhrase.items[0].insertText(\n    \", e na aplica\u00e7\u00e3o dos paradigmas de programa\u00e7\u00e3o \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Edit 2: update the \"Ex-coordenadora\" bullet text. ---\nconst coordPhrase = context.document.body.search(\n  \"Ex-coordenadora da comunidade sem fins lucrativos \",\n  { matchCase: true }\n);\ncoordPhrase.load(\"text\");\nawait context.sync();\n\nif (coordPhrase.items.length > 0) {\n  coordPhrase.items[0].insertText(\n    \"Ex-coordenadora e atual membra da comunidade sem fins lucrativos \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Fix: summary texts highlights\n$d = $word.ActiveDocument\n\n# --- Edit 1: move \"programa\u00e7\u00e3o \" out of the bold run into the preceding\n# non-bold run, leaving only \"funcional\" bold. ---\n\n# Shrink the bold run \"programa\u00e7\u00e3o funcional\" down to just \"funcional\".\n$boldRange = $d.Content\n$boldRange.Find.Execute(\"programa\u00e7\u00e3o funcional\")\n$boldRange.Text = \"funcional\"\n\n# Extend the preceding (non-bold) run so it now ends with \"programa\u00e7\u00e3o \".\n$prefixRange = $d.Content\n$prefixRange.Find.Execute(\", e na aplica\u00e7\u00e3o dos paradigmas de \")\n$prefixRange.Text = \", e na aplica\u00e7\u00e3o dos paradigmas de programa\u00e7\u00e3o \"\n\n# --- Edit 2: update the \"Ex-coordenadora\" bullet text. ---\n$coordRange = $d.Content\n$coordRange.Find.Execute(\"Ex-coordenadora da comunidade sem fins lucrativos \")\n$coordRange.Text = \"Ex-coordenadora e atual membra da comunidade sem fins lucrativos \"\n"}
